$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the data range so that numeric-looking
# strings (e.g. "1.003") are stored as text, matching the original inline-string cells,
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '24.924.01'
$ws.Range('E2').Value = '  +1.15%  '

$ws.Range('D3').Value = '1.705.19'
$ws.Range('E3').Value = '  +0.64%  '

$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.28%  '

$ws.Range('D5').Value = '315.21'
$ws.Range('E5').Value = '  +0.13%  '

$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.26%  '

$ws.Range('D7').Value = '0.4022'
$ws.Range('E7').Value = '  +2.89%  '

$ws.Range('E8').Value = '  +0.55%  '

$ws.Range('D9').Value = '1.006'
$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('D10').Value = '53.79'
$ws.Range('E10').Value = '  +1.42%  '

$ws.Range('D11').Value = '1.471'
$ws.Range('E11').Value = '  -1.79%  '

$ws.Range('D12').Value = '0.08825'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('D13').Value = '26.11'
$ws.Range('E13').Value = '  +6.53%  '

$ws.Range('D14').Value = '7.484'
$ws.Range('E14').Value = '  -2.25%  '

$ws.Range('D15').Value = '8.039'
$ws.Range('E15').Value = '  +0.75%  '

$ws.Range('D16').Value = '0.00001349'
$ws.Range('E16').Value = '  -0.82%  '

$ws.Range('D17').Value = '1.663.43'
$ws.Range('E17').Value = '  -1.70%  '

$ws.Range('D18').Value = '95.43'
$ws.Range('E18').Value = '  -3.16%  '

$ws.Range('D19').Value = '0.07186'
$ws.Range('E19').Value = '  +1.01%  '

$ws.Range('E20').Value = '  +5.98%  '

$ws.Range('D21').Value = '7.264'
$ws.Range('E21').Value = '  -1.02%  '

$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  -0.58%  '

$ws.Range('D23').Value = '14.52'
$ws.Range('E23').Value = '  +1.55%  '

$ws.Range('D24').Value = '24.904.29'
$ws.Range('E24').Value = '  +1.08%  '

$ws.Range('D25').Value = '2.332'
$ws.Range('E25').Value = '  -0.95%  '

$ws.Range('D26').Value = '2.879'
$ws.Range('E26').Value = '  -4.81%  '

$ws.Range('D27').Value = '6.425'
$ws.Range('E27').Value = '  +22.63%  '

$ws.Range('D28').Value = '23.11'
$ws.Range('E28').Value = '  +1.46%  '

$ws.Range('D29').Value = '163.56'
$ws.Range('E29').Value = '  +0.53%  '

$ws.Range('D30').Value = '143.98'
$ws.Range('E30').Value = '  +4.81%  '

$ws.Range('D31').Value = '8.207'
$ws.Range('E31').Value = '  -2.91%  '

$ws.Range('E32').Value = '  +13.82%  '

$ws.Range('D33').Value = '0.08750'
$ws.Range('E33').Value = '  -1.88%  '

$ws.Range('B34').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C34').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D34').Value = '1.830.53'
$ws.Range('E34').Value = '  -2.56%  '

$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '7.374'
$ws.Range('E35').Value = '  -2.23%  '

$ws.Range('D36').Value = '0.03189'
$ws.Range('E36').Value = '  +8.68%  '

$ws.Range('D37').Value = '1.033'
$ws.Range('E37').Value = '  -1.52%  '

$ws.Range('D38').Value = '0.2867'
$ws.Range('E38').Value = '  +4.80%  '

$ws.Range('D39').Value = '0.8530'
$ws.Range('E39').Value = '  +7.87%  '

$ws.Range('D40').Value = '10.85'
$ws.Range('E40').Value = '  +0.39%  '

$ws.Range('D41').Value = '0.09451'
$ws.Range('E41').Value = '  +3.58%  '

$ws.Range('E42').Value = '  -1.31%  '

$ws.Range('D43').Value = '1.475'
$ws.Range('E43').Value = '  +0.56%  '

$ws.Range('D44').Value = '17.82'
$ws.Range('E44').Value = '  +5.41%  '

$ws.Range('D45').Value = '2.720'
$ws.Range('E45').Value = '  +5.64%  '

$ws.Range('D46').Value = '0.7479'
$ws.Range('E46').Value = '  +3.56%  '

$ws.Range('D47').Value = '4.234'
$ws.Range('E47').Value = '  +0.39%  '

$ws.Range('D48').Value = '1.390'
$ws.Range('E48').Value = '  +4.45%  '

$ws.Range('E49').Value = '  -0.29%  '

$ws.Range('D50').Value = '141.19'
$ws.Range('E50').Value = '  +1.55%  '

$ws.Range('D51').Value = '0.08408'
$ws.Range('E51').Value = '  +5.29%  '

# Restore original (unset) cell formatting so the saved file does not carry
# a spurious explicit style on cells that originally had none.
$dataRange.ClearFormats()
